# Update cryptos list D (Price) and E (Volume(1h)) columns to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.401.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.646.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.75%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.646.56"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.131.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.472.93"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.653.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.52"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.778.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "561.10"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.20%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.46%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.22%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0320"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "157.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.95"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0777"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.572"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.92%  "

